$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.449.59'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '1.863.32'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('D5').Value = '310.88'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '0.4774'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.3769'
$ws.Range('E8').Value = '  +2.46%  '
$ws.Range('D9').Value = '0.07329'
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('D10').Value = '0.9355'
$ws.Range('E10').Value = '  +0.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.70'
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '1.895.15'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = '5.433'
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').Value = '6.562'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '90.34'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '0.000008884'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '27.525.12'
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').Value = '14.73'
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').Value = '5.113'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.70'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').Value = '1.944'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').Value = '155.63'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.50'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('D27').Value = '2.025'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = '115.48'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').Value = '4.952'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').Value = '0.08888'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '3.324'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7590'
$ws.Range('E33').Value = '  +2.23%  '
$ws.Range('D34').Value = '4.603'
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('D35').Value = '2.748'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').Value = '0.02059'
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').Value = '0.5635'
$ws.Range('E38').Value = '  +8.15%  '
$ws.Range('D39').Value = '0.05282'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.080'
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').Value = '8.696'
$ws.Range('E42').Value = '  +6.02%  '
$ws.Range('D43').Value = '0.1526'
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '0.4885'
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.67'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '1.661'
$ws.Range('E47').Value = '  +3.40%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '103.24'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').Value = '67.49'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').Value = '0.06077'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').Value = '0.9186'
$ws.Range('E51').Value = '  +3.57%  '
